$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 688.1667
$ws.Range("I12").Value = 176.8
$ws.Range("J12").Value = 1053.4286
$ws.Range("K12").Value = 176.8
$ws.Range("L12").Value = 1053.4286
$ws.Range("M12").Value = -6.800000000000011
$ws.Range("N12").Value = -1393.4286
$ws.Range("H43").Value = 1957.4
$ws.Range("I43").Value = 1300
$ws.Range("J43").Value = 2121.75
$ws.Range("K43").Value = 1300
$ws.Range("L43").Value = 2121.75
$ws.Range("M43").Value = -1231
$ws.Range("N43").Value = -2259.75
$ws.Range("H74").Value = 5609.9546
$ws.Range("I74").Value = 5061.5386
$ws.Range("K74").Value = 5061.5386
$ws.Range("M74").Value = -4125.5386
$ws.Range("H76").Value = 71432160
$ws.Range("I76").Value = 111114360
$ws.Range("J76").Value = 4204
$ws.Range("K76").Value = 111114360
$ws.Range("L76").Value = 4204
$ws.Range("M76").Value = -111114045
$ws.Range("N76").Value = -4834
$ws.Range("H77").Value = 5609.9546
$ws.Range("I77").Value = 5061.5386
$ws.Range("K77").Value = 25307.693
$ws.Range("M77").Value = -20627.693
$ws.Range("H79").Value = 71432160
$ws.Range("I79").Value = 111114360
$ws.Range("J79").Value = 4204
$ws.Range("K79").Value = 111114360
$ws.Range("L79").Value = 4204
$ws.Range("M79").Value = -111113268
$ws.Range("N79").Value = -6388
$ws.Range("H88").Value = 2201.75
$ws.Range("J88").Value = 2201.75
$ws.Range("L88").Value = 2201.75
$ws.Range("N88").Value = -3013.75
$ws.Range("H91").Value = 2201.75
$ws.Range("J91").Value = 2201.75
$ws.Range("L91").Value = 2201.75
$ws.Range("N91").Value = -5009.75
$ws.Range("H132").Value = 3020.1538
$ws.Range("I132").Value = 2890.712
$ws.Range("K132").Value = 8672.136
$ws.Range("M132").Value = -6142.136
$ws.Range("H134").Value = 33994.707
$ws.Range("J134").Value = 33994.707
$ws.Range("L134").Value = 33994.707
$ws.Range("N134").Value = -44134.707
$ws.Range("H138").Value = 4133.9614
$ws.Range("J138").Value = 3966.861
$ws.Range("L138").Value = 11900.583
$ws.Range("N138").Value = -22180.583
$ws.Range("H139").Value = 69999.95
$ws.Range("J139").Value = 69999.95
$ws.Range("L139").Value = 69999.95
$ws.Range("N139").Value = -80279.95

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 9698.929
$ws.Range("I61").Value = 2616.818
$ws.Range("K61").Value = 2616.818
$ws.Range("M61").Value = -2404.818
$ws.Range("H88").Value = 6781.5557
$ws.Range("I88").Value = 1123.1666
$ws.Range("J88").Value = 9610.75
$ws.Range("K88").Value = 1123.1666
$ws.Range("L88").Value = 9610.75
$ws.Range("M88").Value = -717.1666
$ws.Range("N88").Value = -10422.75
$ws.Range("H91").Value = 6781.5557
$ws.Range("I91").Value = 1123.1666
$ws.Range("J91").Value = 9610.75
$ws.Range("K91").Value = 1123.1666
$ws.Range("L91").Value = 9610.75
$ws.Range("M91").Value = 280.8334
$ws.Range("N91").Value = -12418.75
$ws.Range("H136").Value = 9698.929
$ws.Range("I136").Value = 2616.818
$ws.Range("K136").Value = 7850.454000000001
$ws.Range("M136").Value = -5300.454000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3209.5
$ws.Range("I20").Value = 1586.1428
$ws.Range("J20").Value = 6997.3335
$ws.Range("K20").Value = 1586.1428
$ws.Range("L20").Value = 6997.3335
$ws.Range("M20").Value = -1339.1428
$ws.Range("N20").Value = -7491.3335
$ws.Range("H94").Value = 15441.733
$ws.Range("I94").Value = 2108.4
$ws.Range("J94").Value = 22108.4
$ws.Range("K94").Value = 2108.4
$ws.Range("L94").Value = 22108.4
$ws.Range("M94").Value = -1657.4
$ws.Range("N94").Value = -23010.4
$ws.Range("H105").Value = 3781.0908
$ws.Range("I105").Value = 3959.2
$ws.Range("K105").Value = 3959.2
$ws.Range("M105").Value = -2212.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1122.5
$ws.Range("I22").Value = 245.5
$ws.Range("J22").Value = 1999.5
$ws.Range("K22").Value = 245.5
$ws.Range("L22").Value = 1999.5
$ws.Range("M22").Value = 104.5
$ws.Range("N22").Value = -2699.5
$ws.Range("H140").Value = 143592.5
$ws.Range("J140").Value = 143592.5
$ws.Range("L140").Value = 143592.5
$ws.Range("N140").Value = -153952.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 2118999
$ws.Range("J11").Value = 480
$ws.Range("L11").Value = 1440
$ws.Range("N11").Value = -1720
$ws.Range("H38").Value = 191.54546
$ws.Range("J38").Value = 243.42857
$ws.Range("L38").Value = 730.28571
$ws.Range("N38").Value = -1424.28571
$ws.Range("H129").Value = 7577560
$ws.Range("I129").Value = 770.9286
$ws.Range("J129").Value = 20836940
$ws.Range("K129").Value = 2312.7858
$ws.Range("L129").Value = 62510820
$ws.Range("M129").Value = 2687.2142
$ws.Range("N129").Value = -62520820

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 1813.0714
$ws.Range("I80").Value = 1192
$ws.Range("J80").Value = 2278.875
$ws.Range("K80").Value = 1192
$ws.Range("L80").Value = 2278.875
$ws.Range("M80").Value = -194
$ws.Range("N80").Value = -4274.875
$ws.Range("H83").Value = 1813.0714
$ws.Range("I83").Value = 1192
$ws.Range("J83").Value = 2278.875
$ws.Range("K83").Value = 5960
$ws.Range("L83").Value = 11394.375
$ws.Range("M83").Value = -968
$ws.Range("N83").Value = -21378.375
$ws.Range("H102").Value = 3177.9768
$ws.Range("I102").Value = 2784.0356
$ws.Range("K102").Value = 2784.0356
$ws.Range("M102").Value = -1162.0356
$ws.Range("H113").Value = 4779.8
$ws.Range("J113").Value = 4249.75
$ws.Range("L113").Value = 4249.75
$ws.Range("N113").Value = -8589.75

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1521.2142
$ws.Range("J82").Value = 2713
$ws.Range("L82").Value = 2713
$ws.Range("N82").Value = -3435
$ws.Range("H85").Value = 1521.2142
$ws.Range("J85").Value = 2713
$ws.Range("L85").Value = 2713
$ws.Range("N85").Value = -5209
$ws.Range("H93").Value = 1386.5
$ws.Range("J93").Value = 2999
$ws.Range("L93").Value = 2999
$ws.Range("N93").Value = -5495
$ws.Range("H132").Value = 2490.3845
$ws.Range("I132").Value = 2017.5
$ws.Range("K132").Value = 6052.5
$ws.Range("M132").Value = -3522.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 944.5714
$ws.Range("J107").Value = 1209
$ws.Range("L107").Value = 3627
$ws.Range("N107").Value = -7467
$ws.Range("H113").Value = 527365.2
$ws.Range("I113").Value = 770319.4
$ws.Range("J113").Value = 964.3333
$ws.Range("K113").Value = 2310958.2
$ws.Range("L113").Value = 2892.9999
$ws.Range("M113").Value = -2308788.2
$ws.Range("N113").Value = -7232.9999
